$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

$cr = $t.Cell(1, 1).Range
$d.Range($cr.Start, $cr.End - 1).Text = "125÷4=31, 1"
$cr = $t.Cell(1, 2).Range
$d.Range($cr.Start, $cr.End - 1).Text = "367÷8=45, 7"
$cr = $t.Cell(1, 3).Range
$d.Range($cr.Start, $cr.End - 1).Text = "226÷6=37, 4"
$cr = $t.Cell(1, 4).Range
$d.Range($cr.Start, $cr.End - 1).Text = "397÷3=132, 1"
$cr = $t.Cell(1, 5).Range
$d.Range($cr.Start, $cr.End - 1).Text = "364÷7=52, 0"
$cr = $t.Cell(5, 1).Range
$d.Range($cr.Start, $cr.End - 1).Text = "827÷6=137, 5"
$cr = $t.Cell(5, 2).Range
$d.Range($cr.Start, $cr.End - 1).Text = "695÷4=173, 3"
$cr = $t.Cell(5, 3).Range
$d.Range($cr.Start, $cr.End - 1).Text = "351÷7=50, 1"
$cr = $t.Cell(5, 4).Range
$d.Range($cr.Start, $cr.End - 1).Text = "310÷6=51, 4"
$cr = $t.Cell(5, 5).Range
$d.Range($cr.Start, $cr.End - 1).Text = "660÷2=330, 0"
$cr = $t.Cell(9, 1).Range
$d.Range($cr.Start, $cr.End - 1).Text = "950÷2=475, 0"
$cr = $t.Cell(9, 2).Range
$d.Range($cr.Start, $cr.End - 1).Text = "305÷2=152, 1"
$cr = $t.Cell(9, 3).Range
$d.Range($cr.Start, $cr.End - 1).Text = "792÷7=113, 1"
$cr = $t.Cell(9, 4).Range
$d.Range($cr.Start, $cr.End - 1).Text = "912÷5=182, 2"
$cr = $t.Cell(9, 5).Range
$d.Range($cr.Start, $cr.End - 1).Text = "382÷9=42, 4"
$cr = $t.Cell(13, 1).Range
$d.Range($cr.Start, $cr.End - 1).Text = "524÷5=104, 4"
$cr = $t.Cell(13, 2).Range
$d.Range($cr.Start, $cr.End - 1).Text = "588÷7=84, 0"
$cr = $t.Cell(13, 3).Range
$d.Range($cr.Start, $cr.End - 1).Text = "672÷8=84, 0"
$cr = $t.Cell(13, 4).Range
$d.Range($cr.Start, $cr.End - 1).Text = "993÷4=248, 1"
$cr = $t.Cell(13, 5).Range
$d.Range($cr.Start, $cr.End - 1).Text = "231÷3=77, 0"
$cr = $t.Cell(17, 1).Range
$d.Range($cr.Start, $cr.End - 1).Text = "914÷6=152, 2"
$cr = $t.Cell(17, 2).Range
$d.Range($cr.Start, $cr.End - 1).Text = "759÷3=253, 0"
$cr = $t.Cell(17, 3).Range
$d.Range($cr.Start, $cr.End - 1).Text = "742÷9=82, 4"
$cr = $t.Cell(17, 4).Range
$d.Range($cr.Start, $cr.End - 1).Text = "113÷8=14, 1"
$cr = $t.Cell(17, 5).Range
$d.Range($cr.Start, $cr.End - 1).Text = "539÷6=89, 5"
